$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1. Resize the PERT-chart inline picture (2nd inline shape in the
#    document) to the new extent. Unlock the aspect ratio first so
#    width/height can be set independently to match the new image.
# -----------------------------------------------------------------
$pic = $d.InlineShapes.Item(2)
$pic.LockAspectRatio = 0
$pic.Width = 533.374094488189
$pic.Height = 237.79590551181101
$pic.LockAspectRatio = 1

# -----------------------------------------------------------------
# 2. Update the cost paragraph text (new total + cost breakdown).
# -----------------------------------------------------------------
$old = "Τέλος υπολογίζοντας το κόστος του έργου το σύνολο θα κυμανθεί στα 9.500€. Το ποσό αυτό περιλαμβάνει την αμοιβή μας για το διάστημα τον 10,5 μηνών που εκτιμάτε ο χρόνος που χρειαζόμαστε και με τα λειτουργικά κόστη που θα χρειαστούμε."
$new = "Τέλος υπολογίζοντας το κόστος του έργου το σύνολο θα κυμανθεί στα 51.740€. Το ποσό αυτό περιλαμβάνει την αμοιβή μας για το διάστημα τον 10,5 μηνών που εκτιμάτε ο χρόνος που χρειαζόμαστε και με τα λειτουργικά κόστη που θα χρειαστούμε. Το ποσό βγαίνει ως εξής: Υπολογίζονται 1.256 ανθρωποημέρες * 8 ώρες = 10.048 ανθρωποώρες. Θέτοντας ως αμοιβή 5€ την ώρα το ποσό κυμαίνεται στα 50.240€. Τέλος υπολογίζονται 150€ το μήνα τα λειτουργικά έξοδα. Για το διάστημα 10 μηνών τα έξοδα φτάνουν τα 1.500€. Σύνολο οι απολαβές μας θα είναι στα 51.740€."

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
